$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly odds update (FlashScore) - cell value corrections
$updates = @{
    "J2" = 1.13
    "K2" = 6
    "G4" = 3.75
    "H4" = 3.8
    "I4" = 1.82
    "L4" = 1.18
    "M4" = 4.25
    "O4" = 2.27
    "R4" = 1.55
    "S4" = 2.32
    "T4" = 14.5
    "U4" = 24
    "V4" = 12.5
    "W4" = 55
    "X4" = 29
    "Y4" = 29
    "AA4" = 7.6
    "AB4" = 12.5
    "AE4" = 9.75
    "AF4" = 10.5
    "AG4" = 8.25
    "AH4" = 16.5
    "AI4" = 13
    "AJ4" = 19
    "N6" = 1.93
    "O6" = 1.93
    "G10" = 2.1
    "O10" = 1.57
    "G11" = 1.7
    "H11" = 3.4
    "I11" = 4.45
    "N11" = 1.93
    "O11" = 1.7
    "P11" = 1.37
    "Q11" = 2.5
    "R11" = 1.94
    "S11" = 1.77
    "T11" = 5.3
    "U11" = 6.4
    "W11" = 10.75
    "X11" = 11.75
    "Z11" = 8.75
    "AA11" = 5.8
    "AB11" = 13.5
    "AE11" = 9.75
    "AF11" = 20
    "AG11" = 12
    "AH11" = 60
    "AI11" = 35
    "AJ11" = 40
    "G13" = 1.5
    "I13" = 5.7
    "M13" = 3.2
    "N13" = 1.78
    "O13" = 1.82
    "R13" = 2.02
    "S13" = 1.71
    "V13" = 6.8
    "W13" = 9
    "X13" = 10.25
    "Y13" = 21
    "AC13" = 65
    "AE13" = 12
    "AF13" = 28
    "AG13" = 15
    "AH13" = 90
    "AI13" = 50
    "AJ13" = 50
    "N14" = 2.2
    "O14" = 1.65
    "G16" = 2.45
    "W16" = 23
    "X16" = 23
    "K17" = 9
    "N17" = 2.08
    "O17" = 1.73
    "P17" = 1.4
    "Q17" = 2.75
    "Z17" = 9
    "G18" = 2.7
    "I18" = 2.57
    "O18" = 1.47
    "G19" = 2.65
    "I19" = 2.65
    "N19" = 2.05
    "O19" = 1.72
    "G20" = 2.32
    "I20" = 3.1
    "J20" = 1.07
    "K20" = 9
    "N20" = 2.1
    "O20" = 1.67
    "W20" = 23
    "X20" = 21
    "Z20" = 9
    "AD20" = 251
    "AG20" = 11
    "G21" = 2.65
    "H21" = 3.1
    "I21" = 2.65
    "J21" = 1.11
    "K21" = 6.5
    "O21" = 1.47
    "U21" = 12
    "X21" = 26
    "AE21" = 7
    "AF21" = 12
    "AH21" = 29
    "G22" = 2.7
    "H22" = 2.85
    "I22" = 2.65
    "O22" = 1.41
    "G23" = 2.65
    "I23" = 2.3
    "N23" = 1.63
    "I24" = 2.05
    "L24" = 1.33
    "M24" = 3.25
    "N24" = 2.07
    "O24" = 1.69
    "J25" = 1.07
    "L25" = 1.41
    "M25" = 2.62
    "J26" = 1.03
    "L26" = 1.19
    "J27" = 1.04
    "L27" = 1.3
    "J28" = 1.05
    "L28" = 1.41
    "M28" = 2.62
    "J29" = 1.04
    "L29" = 1.27
    "J30" = 1.04
    "L30" = 1.27
    "G31" = 2.05
    "I31" = 3.3
    "J31" = 1.04
    "L31" = 1.3
    "U31" = 9.5
    "X31" = 17
    "AH31" = 41
    "AJ31" = 41
    "G33" = 2.63
    "I33" = 2.55
    "U33" = 15
    "AI33" = 19
    "G34" = 3.5
    "H34" = 3.25
    "I34" = 2.1
    "J34" = 1.1
    "K34" = 6.1
    "L34" = 1.5
    "M34" = 2.5
    "N34" = 2.42
    "O34" = 1.52
    "P34" = 1.52
    "Q34" = 2.42
    "R34" = 2.12
    "S34" = 1.65
    "T34" = 7.9
    "U34" = 17.5
    "V34" = 14
    "W34" = 55
    "X34" = 45
    "Y34" = 65
    "Z34" = 6.1
    "AA34" = 6.7
    "AE34" = 5.7
    "AF34" = 9.25
    "AG34" = 10
    "AH34" = 20
    "AI34" = 23
    "N38" = 1.25
    "N39" = 1.54
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Host "Updated" $updates.Count "cells"